# This workbook tracks weekly Cilantro price observations (Vega Central
# Mapocho de Santiago). A new week of data (two quality-grade rows, one
# priced "$/caja 36 atados" and one priced "$/docena de atados") is added
# at the top of the existing price-history block (which starts at row 230),
# pushing all the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 230-231; every row from the old 230 onward
# shifts down by two (old 230 -> new 232, ... old 335 -> new 337).
$ws.Rows("230:231").Insert()

# --- New row 230: "$/caja 36 atados" observation for the new week ---
$ws.Cells.Item(230, 1).Value  = 9
$ws.Cells.Item(230, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(230, 3).Value  = "Metropolitana"
$ws.Cells.Item(230, 4).Value  = 44466
$ws.Cells.Item(230, 5).Value  = 13
$ws.Cells.Item(230, 6).Value  = 100112040
$ws.Cells.Item(230, 7).Value  = "Cilantro"
$ws.Cells.Item(230, 8).Value  = "Sin especificar"
$ws.Cells.Item(230, 9).Value  = "Primera"
$ws.Cells.Item(230, 10).Value = 52
$ws.Cells.Item(230, 11).Value = 8000
$ws.Cells.Item(230, 12).Value = 8000
$ws.Cells.Item(230, 13).Value = 8000
$ws.Cells.Item(230, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(230, 15).Value = "Región Metropolitana"
$ws.Cells.Item(230, 16).Value = 222
$ws.Cells.Item(230, 17).Value = 36
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# --- New row 231: "$/docena de atados" observation for the new week ---
$ws.Cells.Item(231, 1).Value  = 9
$ws.Cells.Item(231, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(231, 3).Value  = "Metropolitana"
$ws.Cells.Item(231, 4).Value  = 44466
$ws.Cells.Item(231, 5).Value  = 13
$ws.Cells.Item(231, 6).Value  = 100112040
$ws.Cells.Item(231, 7).Value  = "Cilantro"
$ws.Cells.Item(231, 8).Value  = "Sin especificar"
$ws.Cells.Item(231, 9).Value  = "Primera"
$ws.Cells.Item(231, 10).Value = 124
$ws.Cells.Item(231, 11).Value = 10000
$ws.Cells.Item(231, 12).Value = 12000
$ws.Cells.Item(231, 13).Value = 11000
$ws.Cells.Item(231, 14).Value = "$/docena de atados"
$ws.Cells.Item(231, 15).Value = "Región Metropolitana"
$ws.Cells.Item(231, 16).Value = 3667
$ws.Cells.Item(231, 17).Value = 3
$ws.Cells.Item(231, 18).Value = "Hortaliza"
